# "tracks split identifiers column"
#
# The old "Identifiers" column (O) stored compound values like
# "Internal House Name:a round nose". This edit splits that data:
#   - column O is renamed "Internal House Name" and keeps only the
#     value portion (the "Internal House Name:" prefix is dropped)
#   - a new column S "Tag /Band" is appended for the other half of the
#     identifier (left blank for now, to be filled in later)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header and strip the redundant prefix from the data cells.
$ws.Range("O1").Value = "Internal House Name"
$ws.Range("O2").Value = "a round nose"
$ws.Range("O4").Value = "a square nose"
$ws.Range("O5").Value = "a yellow nose"

# New trailing column for the other half of the old identifier.
$ws.Range("S1").Value = "Tag /Band"
